$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D and E) for the new quarters (Dec-2018, Sep-2018).
# This shifts the old D:K data right to F:M, matching Excel native Insert behavior.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting (date / thousands styles) from column F (old D, now shifted)
# into the two freshly inserted columns so D/E match the rest of the table.
$ws.Range("F1:F102").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new Dec-2018 (D) and Sep-2018 (E) quarter data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 834000
$ws.Range("E8").Value = 883000
$ws.Range("D9").Value = 664700
$ws.Range("E9").Value = 698100
$ws.Range("D10").Value = 169300
$ws.Range("E10").Value = 184900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 5300
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 788000
$ws.Range("E17").Value = 812500
$ws.Range("D18").Value = 46000
$ws.Range("E18").Value = 70500
$ws.Range("D20").Value = -30400
$ws.Range("E20").Value = -14900
$ws.Range("D21").Value = 37700
$ws.Range("E21").Value = 77000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 15600
$ws.Range("E23").Value = 55600
$ws.Range("D24").Value = 7300
$ws.Range("E24").Value = 5400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 8300
$ws.Range("E26").Value = 50200
$ws.Range("D27").Value = 8500
$ws.Range("E27").Value = 50200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 2900
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 30400
$ws.Range("E32").Value = 14900
$ws.Range("D33").Value = 11400
$ws.Range("E33").Value = 50200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 11400
$ws.Range("E35").Value = 50200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 170900
$ws.Range("E41").Value = 180800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 413400
$ws.Range("E43").Value = 472400
$ws.Range("D44").Value = 344700
$ws.Range("E44").Value = 334600
$ws.Range("D45").Value = 69800
$ws.Range("E45").Value = 69800
$ws.Range("D46").Value = 998800
$ws.Range("E46").Value = 1057600
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 495400
$ws.Range("E48").Value = 487700
$ws.Range("D49").Value = 1073700
$ws.Range("E49").Value = 1080500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 155400
$ws.Range("E52").Value = 159500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2723300
$ws.Range("E54").Value = 2785300
$ws.Range("D57").Value = 399000
$ws.Range("E57").Value = 419700
$ws.Range("D58").Value = 19400
$ws.Range("E58").Value = 19400
$ws.Range("D59").Value = 139200
$ws.Range("E59").Value = 130700
$ws.Range("D60").Value = 557600
$ws.Range("E60").Value = 569800
$ws.Range("D61").Value = 1336200
$ws.Range("E61").Value = 1316800
$ws.Range("D62").Value = 288900
$ws.Range("E62").Value = 284600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2183300
$ws.Range("E66").Value = 2172000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 472900
$ws.Range("E72").Value = "NA"
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 540000
$ws.Range("E76").Value = 613300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 11400
$ws.Range("E81").Value = 50200
$ws.Range("D83").Value = 22100
$ws.Range("E83").Value = 21400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 82800
$ws.Range("E89").Value = 62000
$ws.Range("D91").Value = -24800
$ws.Range("E91").Value = -19700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -24400
$ws.Range("E94").Value = -15800
$ws.Range("D96").Value = -14000
$ws.Range("E96").Value = -13900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -65100
$ws.Range("E100").Value = -20200
$ws.Range("D101").Value = -3200
$ws.Range("E101").Value = -3800
$ws.Range("D102").Value = -9900
$ws.Range("E102").Value = 22200

# Row 72 (Retained Earnings) source data also marks a few quarters as "NA"
# instead of the previously shifted 0 values, and reports a real number for
# the Dec-2017 quarter (column H).
$ws.Range("F72").Value = "NA"
$ws.Range("G72").Value = "NA"
$ws.Range("H72").Value = 387100
$ws.Range("I72").Value = "NA"
$ws.Range("J72").Value = "NA"
